# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet, and
#    populate it with the fund-holding detail rows (same layout/style as
#    the other quarterly sheets, e.g. "2021-Q4").
# 2. Update the "总计" (Total) summary sheet: insert a new first data row
#    for "2022-Q1" and shift the previous rows down, renumbering column A.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Create the new "2022-Q1" worksheet immediately before "总计".
#
# NOTE: worksheet object references returned by this host are bound to
# a *position*, not a stable identity. Once Worksheets.Add() shifts the
# "总计" sheet from slot 4 to slot 5, any previously-fetched reference to
# it silently starts pointing at whatever now sits in slot 4 (the new
# sheet) instead. So: do the structural change first, and only fetch
# worksheet references to use for writing *after* the sheet collection
# has stopped changing shape.
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($anchor)
$q1.Name = "2022-Q1"

# Re-fetch everything fresh now that sheet order/count is final.
$quarterTemplate = $wb.Worksheets.Item("2021-Q4")   # style/format donor
$totalSheet      = $wb.Worksheets.Item("总计")
$q1              = $wb.Worksheets.Item("2022-Q1")

# Scratch cell used to coerce numeric-looking strings ("29.21", etc.)
# into genuine text values without leaving a quote-prefix style behind:
# write it with a leading apostrophe, copy it, then paste-special
# "values only" into the real destination cell.
function Set-TextValue($sheet, $addr, $val) {
    $sheet.Range("ZZ1").Value = "'" + $val
    $sheet.Range("ZZ1").Copy()
    $sheet.Range($addr).PasteSpecial(-4163)
}

# Header row, copying the formatting used by the other quarter sheets.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
$quarterTemplate.Range("B1:H1").Copy()
$q1.Range("B1").PasteSpecial(-4122)

$fundRows = @(
    @("506002","易方达科创板两年定期开放混合型证券投资基金","29.21","86.29","4.67","1.3641",4),
    @("001513","易方达信息产业混合","32.50","92.37","3.06","0.9945",4),
    @("010013","易方达信息行业精选股票","25.38","88.64","2.49","0.6320",8),
    @("519003","海富通收益增长混合","35.90","78.79","1.53","0.5493",10),
    @("112002","易方达策略成长二号混合","10.72","87.99","4.12","0.4417",4),
    @("110002","易方达策略成长混合","12.15","88.76","3.04","0.3694",9),
    @("501062","南方瑞合三年定期开放混合(LOF)","6.88","56.71","3.18","0.2188",9),
    @("005310","广发电子信息传媒产业精选股票A","3.99","90.16","3.56","0.1420",8),
    @("000328","上投摩根转型动力混合","3.67","82.39","2.84","0.1042",10),
    @("004314","前海开源沪港深新硬件主题灵活配置混合A","1.67","90.05","5.99","0.1000",2),
    @("012200","新华鑫科技3个月滚动持有灵活配置混合型证券投资基金A","2.04","77.02","4.22","0.0861",2),
    @("013757","泰信均衡价值混合A","1.68","59.14","4.51","0.0758",6),
    @("003822","中信建投行业轮换混合A","3.07","72.09","2.36","0.0725",9),
    @("004315","前海开源沪港深新硬件主题灵活配置混合C","1.00","90.05","5.99","0.0599",2),
    @("010236","广发电子信息传媒产业精选股票C","0.81","90.16","3.56","0.0288",8),
    @("007468","中信建投策略精选混合A","0.94","78.13","2.81","0.0264",7),
    @("013758","泰信均衡价值混合C","0.57","59.14","4.51","0.0257",6),
    @("012201","新华鑫科技3个月滚动持有灵活配置混合型证券投资基金C","0.52","77.02","4.22","0.0219",2),
    @("000679","招商丰利灵活配置混合A","0.39","74.75","4.30","0.0168",4),
    @("003823","中信建投行业轮换混合C","0.64","72.09","2.36","0.0151",9),
    @("007469","中信建投策略精选混合C","0.40","78.13","2.81","0.0112",7),
    @("002416","招商丰利灵活配置混合C","0.02","74.75","4.30","0.0009",4)
)

$r = 2
foreach ($row in $fundRows) {
    $q1.Range("A$r").Value = ($r - 2)

    Set-TextValue $q1 "B$r" $row[0]
    $q1.Range("C$r").Value = $row[1]

    Set-TextValue $q1 "D$r" $row[2]
    Set-TextValue $q1 "E$r" $row[3]
    Set-TextValue $q1 "F$r" $row[4]
    Set-TextValue $q1 "G$r" $row[5]

    $q1.Range("H$r").Value = $row[6]

    $r = $r + 1
}

# Copy the "A" column index style (bold/centered) down onto the new rows,
# matching the other quarter sheets' look.
$quarterTemplate.Range("A2").Copy()
$q1.Range("A2:A23").PasteSpecial(-4122)

$q1.Range("ZZ1").Clear()

# ---------------------------------------------------------------------
# Update the "总计" (Total) summary sheet with the new quarter on top.
# ---------------------------------------------------------------------
$totalRows = @(
    @("2022-Q1", 22, 5.36),
    @("2021-Q4", 7, 0.62),
    @("2021-Q3", 5, 0.27),
    @("2021-Q2", 3, 0.21)
)

$r = 2
$idx = 0
foreach ($row in $totalRows) {
    $totalSheet.Range("A$r").Value = $idx
    $totalSheet.Range("B$r").Value = $row[0]
    $totalSheet.Range("C$r").Value = $row[1]
    $totalSheet.Range("D$r").Value = $row[2]
    $r = $r + 1
    $idx = $idx + 1
}
